$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 10): B10:M10 ----
$ws.Range("B10").Value = "Nama"
$ws.Range("C10").Value = "PK"
$ws.Range("D10").Value = "Teknisi 1"
$ws.Range("E10").Value = "Teknisi 2"
$ws.Range("F10").Value = "Teknisi 3"
$ws.Range("G10").Value = "No Hp"
$ws.Range("H10").Value = "Tanggal Perbaikan"
$ws.Range("I10").Value = "Deskripsi Kerusakan"
$ws.Range("J10").Value = "Merk AC"
$ws.Range("K10").Value = "Remot Jenis"
$ws.Range("L10").Value = "Remot Kode"
$ws.Range("M10").Value = "Status"

# copy existing header style (A10) onto the new header cells F10:M10
$ws.Range("A10").Copy()
$ws.Range("F10:M10").PasteSpecial(-4122)

# ---- Row 11 ----
$ws.Range("B11").Value = "Dadan"
$ws.Range("C11").Value = 0.15
$ws.Range("D11").Value = "Abdul Hamid "
$ws.Range("E11").Value = "Erik Hasibuan"
$ws.Range("F11").Value = "Wage Rudolf Supratman "
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "087678987677"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "2025-06-25"
$ws.Range("I11").Value = "Service AC Rutin"
$ws.Range("J11").Value = "Daikin FTC15NV14"
$ws.Range("K11").Value = "Original"
$ws.Range("L11").Value = 27
$ws.Range("M11").Value = 0

# ---- Row 12 ----
$ws.Range("B12").Value = "Muhammad Yamin"
$ws.Range("C12").Value = 0.09
$ws.Range("D12").Value = "Abdul Hamid "
$ws.Range("E12").Value = "Erik Hasibuan"
$ws.Range("F12").Value = "Robert Ed Stewart"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "089976356474"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "2025-06-22"
$ws.Range("I12").Value = "Lorem Ipsum is simply dummy text of the printing and typesetting industry.`n"
$ws.Range("J12").Value = "Changhong CSC-05NVB"
$ws.Range("K12").Value = "Tidak Original"
$ws.Range("L12").Value = 26
$ws.Range("M12").Value = 0

# ---- Row 13 ----
$ws.Range("B13").Value = "Jihan Fahriza Amalina "
$ws.Range("C13").Value = 0.09
$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = "-"
$ws.Range("F13").Value = "Abdul Yamin"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "087898876567"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "2025-06-18"
$ws.Range("I13").Value = "Ganti Preon"
$ws.Range("J13").Value = "Denpoo DDS-199CI"
$ws.Range("K13").Value = "Original"
$ws.Range("L13").Value = 25
$ws.Range("M13").Value = 0

# ---- Row 14 ----
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "Usep"
$ws.Range("C14").Value = 0.12
$ws.Range("D14").Value = "Qarib Abdullah Shakil"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "Wage Rudolf Supratman "
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "08976756765"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "2025-06-20"
$ws.Range("I14").Value = "Lorem Ipsum is simply dummy text of the printing and typesetting industry."
$ws.Range("J14").Value = "LG H05TN4"
$ws.Range("K14").Value = "Original"
$ws.Range("L14").Value = 24
$ws.Range("M14").Value = 1

# copy data-row style (E11, which already has s=5) onto all newly written data cells
$ws.Range("E11").Copy()
$ws.Range("F11:M14").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("A14").PasteSpecial(-4122)

# move selection to M14 to match the final cursor position
$ws.Range("M14").Select()
